# Updates cryptos list prices (column D) and 1h volume change (column E)
# to match the latest scrape, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(NewPrice-or-$null, NewVolumePercent)
# NewPrice is $null when column D did not change for that row.
$updates = @{
    2  = @("61.331.11", "  +0.62%  ")
    3  = @("2.931.89",  "  +0.51%  ")
    4  = @($null,       "  -0.04%  ")
    5  = @("594.82",    "  +0.36%  ")
    6  = @("143.54",    "  -1.44%  ")
    7  = @($null,       "  +0.01%  ")
    8  = @($null,       "  -1.16%  ")
    9  = @($null,       "  +0.95%  ")
    10 = @($null,       "  -1.77%  ")
    11 = @($null,       "  -0.51%  ")
    12 = @($null,       "  -0.97%  ")
    13 = @("33.24",     "  -0.94%  ")
    14 = @($null,       "  +0.29%  ")
    15 = @("3.416.69",  "  +0.50%  ")
    16 = @("61.318.89", "  +0.61%  ")
    17 = @("2.930.99",  "  +0.48%  ")
    18 = @($null,       "  -0.59%  ")
    19 = @($null,       "  +0.79%  ")
    20 = @("13.56",     "  +1.50%  ")
    21 = @($null,       "  -1.08%  ")
    22 = @($null,       "  -0.11%  ")
    23 = @("81.45",     "  +0.09%  ")
    24 = @("10.86",     "  -0.60%  ")
    25 = @("2.17",      "  -1.93%  ")
    26 = @($null,       "  -1.96%  ")
    27 = @($null,       "  +0.00%  ")
    28 = @("2.20",      "  -3.85%  ")
    29 = @($null,       "  -0.89%  ")
    30 = @($null,       "  -2.40%  ")
    31 = @("26.64",     "  +0.82%  ")
    32 = @($null,       "  +1.32%  ")
    33 = @($null,       "  +0.01%  ")
    34 = @("0.0₃0872",  "  +2.11%  ")
    35 = @($null,       "  -0.36%  ")
    36 = @("5.62",      "  -0.23%  ")
    37 = @($null,       "  -2.10%  ")
    38 = @($null,       "  -0.15%  ")
    39 = @($null,       "  +0.34%  ")
    40 = @("8.50",      "  -0.52%  ")
    41 = @("42.14",     "  +5.54%  ")
    42 = @("0.278",     "  -3.10%  ")
    43 = @("2.703.72",  "  +0.10%  ")
    44 = @($null,       "  -0.40%  ")
    45 = @("133.46",    "  +1.21%  ")
    46 = @("363.53",    "  -3.11%  ")
    48 = @("23.52",     "  -1.48%  ")
    49 = @($null,       "  -1.29%  ")
    50 = @($null,       "  -0.50%  ")
    51 = @($null,       "  +0.01%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        # Prefix with an apostrophe so Excel keeps the value as literal text
        # (matching the original inline-string/text storage) instead of
        # silently re-interpreting numeric-looking strings such as
        # "594.82" as a floating point number.
        $ws.Cells.Item($row, 4).Value = "'" + $newPrice
    }
    $ws.Cells.Item($row, 5).Value = $newVolume
}
